$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Add the new row of log data (row 5)
$logs.Range("A5").Value = "Kun je dit voor me fixen?"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #1: Kun je dit voor me fixen?"
$logs.Range("D5").Value = "Planning / Afspraak"
$logs.Range("E5").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Range("F5").Value = "2025-08-06 19:36:15"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Ja"
$logs.Range("I5").Value = "Nee"
$logs.Range("J5").Value = "Nee"

# Extend the existing conditional formatting ranges to include the new row
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "4")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "5")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count
$dashboard.Range("B2").Value = 4
